# Generate Report for Handoff
# The handback for e2e\f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md is now ready
# for handoff again (its handed-back version is stale vs. the latest source),
# so update its status + timestamps on all three sheets and record the
# "stale handback" error detail on the per-language sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/8fbaa06a4c6077b31f4b47282fbc5a9429f70e45/e2e/f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e91464516753726870e0065f112be94a88620c58/e2e/f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md."

# --- Overview sheet: row 3 is the f5fcffff-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 21:03:18"

# --- zh-cn sheet: row 3 is the f5fcffff-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-13 21:03:10"
$wsZhCn.Range("P3").Value = $errorDetail
# ColumnWidth uses Excel's "characters" scale, which renders as a stored
# sheet width of 40 (the same scale/units used by <col width="...">) when
# set to 39.17 for the workbook's default font (Calibri 11).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the f5fcffff-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-13 21:03:18"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
